# Apply "update to supplemental tables" edits.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Table S1 - Plasticity AIC (sheet 1)
#   - model formula text gains " + (1 | tank)" random-effect term
#   - column B widened to fit the longer formula text
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Table S1 - Plasticity AIC")
$ws1.Range("B2").Value = "reef environment * pCO2 + temperature + (1 | colony) + (1 | tank)"
$ws1.Range("B16").Value = "pCO2 + temperature + (1 | colony) + (1 | tank)"
# Target raw column width is 65.71 characters. The host's ColumnWidth setter
# snaps to a 6px/char + 5px pixel grid, so 64.8 is the closest settable value
# that lands on the stored width nearest to 65.71 (65.67).
$ws1.Columns.Item(2).ColumnWidth = 64.8

# ---------------------------------------------------------------------
# Table S2 - PERMANOVA (sheet 2) - updated statistics
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Table S2 - PERMANOVA")

$ws2.Range("C2").Value = 59423
$ws2.Range("D2").Value = 0.203
$ws2.Range("E2").Value = 7.93

$ws2.Range("C3").Value = 9320
$ws2.Range("D3").Value = 0.032
$ws2.Range("E3").Value = 3.73
$ws2.Range("F3").Value = 0.06529

$ws2.Range("C7").Value = 101796
$ws2.Range("D7").Value = 0.09
$ws2.Range("E7").Value = 14.87

$ws2.Range("C8").Value = 519372
$ws2.Range("D8").Value = 0.46
$ws2.Range("E8").Value = 75.84

$ws2.Range("C12").Value = 724
$ws2.Range("D12").Value = 0.005
$ws2.Range("E12").Value = 0.53
$ws2.Range("F12").Value = 0.45237

$ws2.Range("C13").Value = 27051
$ws2.Range("D13").Value = 0.191
$ws2.Range("E13").Value = 19.66

# ---------------------------------------------------------------------
# Table S4 - Species PERMANOVA (sheet 4) - updated p-values
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Table S4 - Species PERMANOVA")

$ws4.Range("F3").Value = 0.09327
$ws4.Range("F4").Value = 0.00466
$ws4.Range("F7").Value = 0.02398

# ---------------------------------------------------------------------
# Table S5 - HostVsymb PERMANOVA (sheet 5) - updated p-values
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Table S5 - HostVsymb PERMANOVA")

$ws5.Range("F2").Value = 0.74883

$ws5.Range("F3").Value = 0.00333
$ws5.Range("K3").Value = 0.09927

$ws5.Range("F4").Value = 0.57295
$ws5.Range("K4").Value = 0.002

$ws5.Range("F7").Value = 0.30713
$ws5.Range("K7").Value = 0.27981

$ws5.Range("F9").Value = 0.15723

$ws5.Range("F12").Value = 0.01599
$ws5.Range("K12").Value = 0.00133

$ws5.Range("F13").Value = 0.08328

$ws5.Range("F14").Value = 0.18055
$ws5.Range("K14").Value = 0.47901

Write-Host "Supplemental tables updated"
